$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.191405653953552
$ws.Range("B1").Value = 1.723453044891357
$ws.Range("C1").Value = 6.817564010620117
$ws.Range("D1").Value = 2.273456573486328
$ws.Range("E1").Value = 1.190638542175293
